$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates ---
$ws.Range("G5").Value = 61.34

$ws.Range("B13").Value = 4433450
$ws.Range("C13").Value = 9604214

$ws.Range("B14").Value = 9724094
$ws.Range("C14").Value = 8514602

$ws.Range("B15").Value = 6974566
$ws.Range("C15").Value = 1113865

$ws.Range("B16").Value = 7287140
$ws.Range("C16").Value = 6927695

$ws.Range("E12").Value = "Example Narrative:"
$ws.Range("E14").Value = "sample narrative for the report"

Write-Output "done"
